# (NCIOCPL#124) Fixed data errors for citations and language toggle.
#
# - Corrects the Spanish "expectedHeaderText" values on the
#   pages_with_citations sheet: the Spanish article row should read
#   "Bibliografía selecta" (was incorrectly "Selected References") and the
#   Spanish press-release row should read "Bibliografía" (was the
#   untranslated "Referencias").
# - Restores the view/selection state that Excel persisted when the file
#   was last saved (active sheet/tab and selected ranges).

$wb = $excel.ActiveWorkbook

$wsCitations = $wb.Worksheets.Item("pages_with_citations")
$wsNoCitations = $wb.Worksheets.Item("pages_without_citations")

# Fix the mistranslated / mismatched expected header text values.
$wsCitations.Range("D4").Value = "Bibliografía selecta"
$wsCitations.Range("D5").Value = "Bibliografía"

# Restore sheet view/selection state: pages_without_citations was last
# active with C12 selected, then pages_with_citations becomes the active
# (selected) tab with D4:D5 selected.
$wsNoCitations.Activate() | Out-Null
$wsNoCitations.Range("C12").Select() | Out-Null

$wsCitations.Activate() | Out-Null
$wsCitations.Range("D4:D5").Select() | Out-Null
